# Shift the "arma_ns_primeira_letra" column (A) values forward by 3 letters
# of the alphabet for rows 350 through 469 (A->D, B->E, C->F, ... J->M).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 350; $row -le 469; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $current = [string]$cell.Value2
    if ($current -ne "") {
        $code = [int][char]$current
        $newChar = [char]($code + 3)
        $cell.Value = [string]$newChar
    }
}
